$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (pushes old rows 5-12 down to 6-13,
# carrying their values/styles/heights along unchanged).
$null = $ws.Rows.Item(5).Insert()

# Row 4: task description / estimate changed (dependency fix)
$ws.Range("B4").Value = "Inicijalizacija projekta za pristup bazi podataka, osnovna podesavanja, implementacija generickog repozitorijuma, UnitOfWorka, DbContexta"
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 35

# New row 5: newly added task
$ws.Range("B5").Value = "Implementacija repozitorijuma za Wallet i Tranasaction, dodavanje Tabela i konfiguracija za Wallet i Transaction"
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 15
$ws.Rows.Item(5).RowHeight = 31.5

# Update selection to match the saved workbook view
$null = $ws.Range("B6").Select()
